$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.181.15'
$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").Value = '1.640.29'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.523'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.22%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").Value = '1.870.42'
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").Value = '1.644.47'
$ws.Range("E13").Value = '  +0.07%  '

$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").Value = '27.201.26'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("E24").Value = '  -0.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("E29").Value = '  -0.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0509'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.98%  '

$ws.Range("E31").Value = '  +0.13%  '

$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.28%  '

$ws.Range("D34").Value = '1.308.35'
$ws.Range("E34").Value = '  +3.63%  '

$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("E36").Value = '  +1.51%  '

$ws.Range("E37").Value = '  -1.40%  '

$ws.Range("E38").Value = '  +2.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.546'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.95%  '

$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.67%  '

$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").Value = '1.780.52'
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '

$ws.Range("E47").Value = '  +1.29%  '

$ws.Range("E48").Value = '  +1.12%  '

$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0963'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.19%  '
